$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-08-12 Monday" "2024-08-13 Tuesday"

Replace-Text "850×8=" "832×5="
Replace-Text "145×6=" "556×4="
Replace-Text "737×4=" "231×9="
Replace-Text "949×2=" "462×3="
Replace-Text "964×3=" "552×4="

Replace-Text "927×6=" "548×6="
Replace-Text "358×7=" "943×7="
Replace-Text "537×2=" "817×6="
Replace-Text "154×4=" "577×8="
Replace-Text "281×7=" "182×3="

Replace-Text "975×2=" "107×9="
Replace-Text "997×9=" "899×2="
Replace-Text "363×6=" "301×3="
Replace-Text "631×3=" "453×3="
Replace-Text "578×9=" "834×4="

Replace-Text "626×4=" "534×7="
Replace-Text "920×2=" "249×3="
Replace-Text "152×3=" "158×5="
Replace-Text "976×8=" "765×9="
Replace-Text "864×4=" "870×6="

Replace-Text "285×9=" "505×9="
Replace-Text "376×5=" "993×5="
Replace-Text "270×6=" "927×9="
Replace-Text "296×4=" "368×8="
Replace-Text "878×8=" "389×8="
